# Normalize the "Recorded By" (column G) cell values: the comma-separated
# list of recorders/emails in each row gets sorted alphabetically
# (case-insensitive), e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Text

    if ([string]::IsNullOrEmpty($current)) {
        continue
    }

    $parts = $current -split ", "
    $sorted = $parts | Sort-Object
    $newValue = $sorted -join ", "

    if ($newValue -ne $current) {
        $cell.Value = $newValue
    }
}

Write-Host "Recorded By column normalized through row $lastRow"
